# Auto-generated Excel COM-interop script applying the Ravana_Profits value updates
# produced by the scheduled runner. Sets explicit numeric values per cell, and clears
# the handful of cells whose column no longer applies for that row (matching the diff).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 250
$ws.Cells.Item(18, 9).Value = 250
$ws.Cells.Item(18, 11).Value = 250
$ws.Cells.Item(18, 13).Value = 34
$ws.Cells.Item(64, 8).Value = 3700
$ws.Cells.Item(64, 10).Value = 3700
$ws.Cells.Item(64, 12).Value = 3700
$ws.Cells.Item(64, 14).Value = -4196
$ws.Cells.Item(67, 8).Value = 3700
$ws.Cells.Item(67, 10).Value = 3700
$ws.Cells.Item(67, 12).Value = 3700
$ws.Cells.Item(67, 14).Value = -5416
$ws.Cells.Item(113, 8).Value = 3428.4285
$ws.Cells.Item(113, 10).Value = 3416.6667
$ws.Cells.Item(113, 12).Value = 3416.6667
$ws.Cells.Item(113, 14).Value = -9924.6667
$ws.Cells.Item(129, 8).Value = 2574
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(129, 14).ClearContents()
$ws.Cells.Item(138, 8).Value = 4533.4
$ws.Cells.Item(138, 9).Value = 3170.1428
$ws.Cells.Item(138, 10).Value = 5063.5557
$ws.Cells.Item(138, 11).Value = 9510.428400000001
$ws.Cells.Item(138, 12).Value = 15190.6671
$ws.Cells.Item(138, 13).Value = -4370.428400000001
$ws.Cells.Item(138, 14).Value = -25470.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 420
$ws.Cells.Item(5, 9).Value = 420
$ws.Cells.Item(5, 11).Value = 420
$ws.Cells.Item(5, 13).Value = -308
$ws.Cells.Item(32, 8).Value = 17004.074
$ws.Cells.Item(32, 9).Value = 16764.48
$ws.Cells.Item(32, 11).Value = 16764.48
$ws.Cells.Item(32, 13).Value = -16477.48
$ws.Cells.Item(74, 8).Value = 1697.1724
$ws.Cells.Item(74, 9).Value = 1134.2727
$ws.Cells.Item(74, 11).Value = 1134.2727
$ws.Cells.Item(74, 13).Value = -260.2727
$ws.Cells.Item(77, 8).Value = 1697.1724
$ws.Cells.Item(77, 9).Value = 1134.2727
$ws.Cells.Item(77, 11).Value = 5671.363499999999
$ws.Cells.Item(77, 13).Value = -1303.363499999999
$ws.Cells.Item(132, 8).Value = 4221.5
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 4221.5
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 12664.5
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -17724.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 420
$ws.Cells.Item(4, 9).Value = 420
$ws.Cells.Item(4, 11).Value = 420
$ws.Cells.Item(4, 13).Value = -305
$ws.Cells.Item(22, 8).Value = 431.66666
$ws.Cells.Item(22, 9).Value = 431.66666
$ws.Cells.Item(22, 11).Value = 431.66666
$ws.Cells.Item(22, 13).Value = -258.66666
$ws.Cells.Item(64, 8).Value = 1164.1666
$ws.Cells.Item(64, 10).Value = 1326.3334
$ws.Cells.Item(64, 12).Value = 1326.3334
$ws.Cells.Item(64, 14).Value = -1776.3334
$ws.Cells.Item(67, 8).Value = 1164.1666
$ws.Cells.Item(67, 10).Value = 1326.3334
$ws.Cells.Item(67, 12).Value = 1326.3334
$ws.Cells.Item(67, 14).Value = -2886.3334
$ws.Cells.Item(97, 8).Value = 10000
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 10000
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 10000
$ws.Cells.Item(97, 13).ClearContents()
$ws.Cells.Item(97, 14).Value = -11982
$ws.Cells.Item(99, 8).Value = 1499.5
$ws.Cells.Item(99, 9).Value = 1000
$ws.Cells.Item(99, 11).Value = 1000
$ws.Cells.Item(99, 13).Value = 498
$ws.Cells.Item(105, 8).Value = 7670.857
$ws.Cells.Item(105, 9).Value = 7616
$ws.Cells.Item(105, 11).Value = 7616
$ws.Cells.Item(105, 13).Value = -5869
$ws.Cells.Item(137, 8).Value = 99999.5
$ws.Cells.Item(137, 10).Value = 99999.5
$ws.Cells.Item(137, 12).Value = 99999.5
$ws.Cells.Item(137, 14).Value = -110199.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 650
$ws.Cells.Item(22, 9).Value = 650
$ws.Cells.Item(22, 11).Value = 650
$ws.Cells.Item(22, 13).Value = -300
$ws.Cells.Item(31, 8).Value = 2026.037
$ws.Cells.Item(31, 9).Value = 1988.6538
$ws.Cells.Item(31, 11).Value = 1988.6538
$ws.Cells.Item(31, 13).Value = -1693.6538
$ws.Cells.Item(34, 8).Value = 2026.037
$ws.Cells.Item(34, 9).Value = 1988.6538
$ws.Cells.Item(34, 11).Value = 1988.6538
$ws.Cells.Item(34, 13).Value = -1786.6538
$ws.Cells.Item(53, 8).Value = 61842
$ws.Cells.Item(53, 10).Value = 61842
$ws.Cells.Item(53, 12).Value = 61842
$ws.Cells.Item(53, 14).Value = -63056
$ws.Cells.Item(62, 8).Value = 4493
$ws.Cells.Item(62, 10).Value = 4336
$ws.Cells.Item(62, 12).Value = 4336
$ws.Cells.Item(62, 14).Value = -5584
$ws.Cells.Item(65, 8).Value = 4493
$ws.Cells.Item(65, 10).Value = 4336
$ws.Cells.Item(65, 12).Value = 21680
$ws.Cells.Item(65, 14).Value = -27920
$ws.Cells.Item(86, 8).Value = 23358.947
$ws.Cells.Item(86, 9).Value = 10989.777
$ws.Cells.Item(86, 10).Value = 34491.2
$ws.Cells.Item(86, 11).Value = 10989.777
$ws.Cells.Item(86, 12).Value = 34491.2
$ws.Cells.Item(86, 13).Value = -9866.777
$ws.Cells.Item(86, 14).Value = -36737.2
$ws.Cells.Item(89, 8).Value = 23358.947
$ws.Cells.Item(89, 9).Value = 10989.777
$ws.Cells.Item(89, 10).Value = 34491.2
$ws.Cells.Item(89, 11).Value = 54948.885
$ws.Cells.Item(89, 12).Value = 172456
$ws.Cells.Item(89, 13).Value = -49332.885
$ws.Cells.Item(89, 14).Value = -183688
$ws.Cells.Item(99, 8).Value = 9473.666999999999
$ws.Cells.Item(99, 9).Value = 9281
$ws.Cells.Item(99, 11).Value = 9281
$ws.Cells.Item(99, 13).Value = -7783
$ws.Cells.Item(126, 8).Value = 9473.666999999999
$ws.Cells.Item(126, 9).Value = 9281
$ws.Cells.Item(126, 11).Value = 27843
$ws.Cells.Item(126, 13).Value = -25373
$ws.Cells.Item(141, 8).Value = 60710.465
$ws.Cells.Item(141, 10).Value = 60710.465
$ws.Cells.Item(141, 12).Value = 60710.465
$ws.Cells.Item(141, 14).Value = -71070.465

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 699.1667
$ws.Cells.Item(5, 9).Value = 831.6667
$ws.Cells.Item(5, 10).Value = 566.6667
$ws.Cells.Item(5, 11).Value = 2495.0001
$ws.Cells.Item(5, 12).Value = 1700.0001
$ws.Cells.Item(5, 13).Value = -2383.0001
$ws.Cells.Item(5, 14).Value = -1924.0001
$ws.Cells.Item(107, 8).Value = 2230.4
$ws.Cells.Item(107, 9).Value = 5101.5
$ws.Cells.Item(107, 10).Value = 316.33334
$ws.Cells.Item(107, 11).Value = 15304.5
$ws.Cells.Item(107, 12).Value = 949.0000200000001
$ws.Cells.Item(107, 13).Value = -13384.5
$ws.Cells.Item(107, 14).Value = -4789.00002
$ws.Cells.Item(113, 8).Value = 4067.8572
$ws.Cells.Item(113, 10).Value = 4080
$ws.Cells.Item(113, 12).Value = 12240
$ws.Cells.Item(113, 14).Value = -16580
$ws.Cells.Item(135, 8).Value = 699.1667
$ws.Cells.Item(135, 9).Value = 831.6667
$ws.Cells.Item(135, 10).Value = 566.6667
$ws.Cells.Item(135, 11).Value = 7485.0003
$ws.Cells.Item(135, 12).Value = 5100.0003
$ws.Cells.Item(135, 13).Value = -4950.0003
$ws.Cells.Item(135, 14).Value = -10170.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 713.9231
$ws.Cells.Item(97, 9).Value = 663.2
$ws.Cells.Item(97, 10).Value = 883
$ws.Cells.Item(97, 11).Value = 663.2
$ws.Cells.Item(97, 12).Value = 883
$ws.Cells.Item(97, 13).Value = -167.2
$ws.Cells.Item(97, 14).Value = -1875
$ws.Cells.Item(122, 8).Value = 2003.5
$ws.Cells.Item(122, 9).Value = 2003.5
$ws.Cells.Item(122, 11).Value = 6010.5
$ws.Cells.Item(122, 13).Value = -3560.5
$ws.Cells.Item(126, 8).Value = 5375.75
$ws.Cells.Item(126, 9).Value = 5302.75
$ws.Cells.Item(126, 10).Value = 5448.75
$ws.Cells.Item(126, 11).Value = 15908.25
$ws.Cells.Item(126, 12).Value = 16346.25
$ws.Cells.Item(126, 13).Value = -13438.25
$ws.Cells.Item(126, 14).Value = -21286.25
$ws.Cells.Item(132, 8).Value = 3626.3333
$ws.Cells.Item(132, 9).Value = 2253.6667
$ws.Cells.Item(132, 10).Value = 4999
$ws.Cells.Item(132, 11).Value = 6761.000100000001
$ws.Cells.Item(132, 12).Value = 14997
$ws.Cells.Item(132, 13).Value = -4231.000100000001
$ws.Cells.Item(132, 14).Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 448.66666
$ws.Cells.Item(46, 9).Value = 448.66666
$ws.Cells.Item(46, 11).Value = 448.66666
$ws.Cells.Item(46, 13).Value = -260.66666
$ws.Cells.Item(61, 8).Value = 4014.625
$ws.Cells.Item(61, 9).Value = 3945.2856
$ws.Cells.Item(61, 10).Value = 4500
$ws.Cells.Item(61, 11).Value = 3945.2856
$ws.Cells.Item(61, 12).Value = 4500
$ws.Cells.Item(61, 13).Value = -3743.2856
$ws.Cells.Item(61, 14).Value = -4904
$ws.Cells.Item(113, 8).Value = 4014.625
$ws.Cells.Item(113, 9).Value = 3945.2856
$ws.Cells.Item(113, 10).Value = 4500
$ws.Cells.Item(113, 11).Value = 3945.2856
$ws.Cells.Item(113, 12).Value = 4500
$ws.Cells.Item(113, 13).Value = -1775.2856
$ws.Cells.Item(113, 14).Value = -8840

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3070.3125
$ws.Cells.Item(132, 9).Value = 1963.9
$ws.Cells.Item(132, 10).Value = 4914.3335
$ws.Cells.Item(132, 11).Value = 5891.700000000001
$ws.Cells.Item(132, 12).Value = 14743.0005
$ws.Cells.Item(132, 13).Value = -3361.700000000001
$ws.Cells.Item(132, 14).Value = -19803.0005

